$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows (bottom-up to avoid shifting issues) for Caso 5828 (row19), Caso 5642 (row15), Caso 3289 (row9)
$ws.Rows.Item(19).Delete()
$ws.Rows.Item(15).Delete()
$ws.Rows.Item(9).Delete()
